$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: "Send a calendar file" ---
# Progress status changed from a text status ("In Progress") to a 60% completion
# figure, and a note was added describing current state.
$ws.Range("E50").NumberFormat = "0%"
$ws.Range("E50").Value = 0.6
$ws.Range("F50").Value = "Code compiles; waiting on SMTP IP"

# --- Row 54: "What data do we need from user / prof to send receive invites" ---
$ws.Range("E54").Value = 0.6
$ws.Range("F54").Value = "Class has most data"

# --- Row 51: "Validate working with Gannon outlook" ---
$ws.Range("F51").Value = "Anticipated, Waiting on SMPT IP"

# --- Row 52: "Validate accept/decline uses" ---
$ws.Range("F52").Value = "Anticipated, waiting on SMTP IP"

# --- Row 53: "Define best user flow for sending receiving" ---
$ws.Range("E53").Value = 0.6

# Update the active selection to reflect where the edits were made.
$ws.Range("F53").Select()
